# Applies the 2022-12-25 symbol-list refresh to Sheet1 (cryptos.xlsx).
# Values in this sheet are stored as text (inline strings), including
# numeric-looking prices, so every write below goes through a helper
# that keeps numeric-looking text as TEXT (Excel would otherwise
# auto-convert a plain numeric string into a Number cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($Cell, $Text) {
    $range = $ws.Range($Cell)
    if ($Text -match '^-?[0-9]+(\.[0-9]+)?$') {
        # Force text storage with a leading apostrophe, then strip the
        # resulting quote-prefix formatting so the cell style is untouched.
        $range.Value = "'" + $Text
        $range.ClearFormats()
    } else {
        $range.Value = $Text
    }
}

Set-CellText 'D3' '23.04'
Set-CellText 'D4' '5.410'
Set-CellText 'D5' '0.06022'
Set-CellText 'D6' '3.393'
Set-CellText 'D7' '0.8109'
Set-CellText 'D8' '0.9288'
Set-CellText 'D11' '0.03369'
Set-CellText 'D12' '0.03038'
Set-CellText 'D13' '0.09352'
Set-CellText 'D15' '0.001593'
Set-CellText 'D40' '0.03967'
Set-CellText 'D41' '0.006413'
Set-CellText 'D42' '0.1072'
Set-CellText 'D43' '0.002901'
Set-CellText 'D44' '0.006558'
Set-CellText 'E44' '43LocalTradersLCT'
Set-CellText 'D45' '0.00005206'
Set-CellText 'D47' '0.0005802'
Set-CellText 'D48' '0.8503'
Set-CellText 'D49' '0.002291'
Set-CellText 'B17' 'TigerCash'
Set-CellText 'C17' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-CellText 'D17' '0.005372'
Set-CellText 'E17' '16TigerCashTCH'
Set-CellText 'B18' 'HotbitToken'
Set-CellText 'C18' 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-CellText 'D18' '0.004151'
Set-CellText 'E18' '17HotbitTokenHTB'
Set-CellText 'B19' 'BitKan'
Set-CellText 'C19' 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-CellText 'D19' '0.0009826'
Set-CellText 'E19' '18BitKanKAN'
Set-CellText 'B20' 'NitroEx'
Set-CellText 'C20' 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-CellText 'D20' '0.00008703'
Set-CellText 'E20' '19NitroExNTX'
Set-CellText 'B21' 'LEO'
Set-CellText 'C21' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-CellText 'D21' '3.655'
Set-CellText 'E21' '20LEOLEO'
Set-CellText 'B22' 'KuCoinToken'
Set-CellText 'C22' 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-CellText 'D22' '6.442'
Set-CellText 'E22' '21KuCoinTokenKCS'
Set-CellText 'B23' 'BTSEToken'
Set-CellText 'C23' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-CellText 'D23' '2.185'
Set-CellText 'E23' '22BTSETokenBTSE'
Set-CellText 'B24' 'One'
Set-CellText 'C24' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-CellText 'D24' '0.01118'
Set-CellText 'E24' '23OneONEBestin24h'
